# Update automatico via Actualizar 02-06-2021 02-36-51
#
# Mirrors the periodic "disponibilidad" scrape: the last 14-row availability
# block (rows 940-953, a re-check of the same 14 services) gets a refreshed
# timestamp (tiny float re-serialization) and a brand-new 14-row block
# (rows 954-967) is appended for the newest check, cycling through the same
# services/URLs as every previous block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 940:953 - D column (Fecha) gets re-stamped with a slightly more
#    precise float serialization of the same instant.
# ---------------------------------------------------------------------
for ($r = 940; $r -le 953; $r++) {
    $ws.Range("D$r").Value2 = 44233.08765997685
}

# ---------------------------------------------------------------------
# 2) Append rows 954:967 - same 14-service cycle as every earlier block
#    (see rows 2:15, 16:29, ...): column A = service name (shared string),
#    column B = URL (hyperlinked, shared string), column C = "Disponible",
#    column D = new check timestamp.
# ---------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Hyperlink "Address" target: identical to $urls, except the MapStore entry
# is split into a bare address + a "/" sub-address/location (matches how
# every earlier occurrence of that same link was stored).
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkSubAddresses = @($null,$null,$null,$null,$null,$null,$null,$null,"/",$null,$null,$null,$null,$null)

$newTimestamp = 44233.10887606165
$startRow = 954

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value2 = $names[$i]
    $ws.Range("C$r").Value2 = "Disponible"

    $dCell = $ws.Range("D$r")
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dCell.Value2 = $newTimestamp

    $bCell = $ws.Range("B$r")
    # Pre-set the display text to the full URL (incl. any #fragment) so the
    # shared-string the hyperlink sits over matches the one reused by every
    # earlier row with this same link - Hyperlinks.Add() below must not be
    # left to invent its own display text.
    $bCell.Value2 = $urls[$i]

    $sub = $linkSubAddresses[$i]
    if ($sub) {
        $h = $ws.Hyperlinks.Add($bCell, $linkAddresses[$i], $sub)
    } else {
        $h = $ws.Hyperlinks.Add($bCell, $linkAddresses[$i])
    }

    # Hyperlinks.Add() always reformats the cell with a freshly minted style;
    # snap it back onto the workbook's single pre-existing "Hyperlink" cell
    # format (style index used by every B2:B953 link) instead of leaving the
    # new one in place.
    $bCell.Style = "Hyperlink"
}
